$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$h = $ws.Hyperlinks
$item = $h.Item(9)
$r = $item.Range()
"addr=$($r.Address())"
$item.Delete()
"count after delete: $($h.Count())"
